$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bibliography rows appended below the existing 3 rows of data.
$ws.Range("A4").Value = "An Industrial Application of Mutation Testing - Lessons, Challenges, and Research Directions"
$ws.Range("B4").Value = 'G. Petrovic, M. Ivankovic, B. Kurtz, P. Ammann and R. Just, "An Industrial Application of Mutation Testing: Lessons, Challenges, and Research Directions," 2018 IEEE International Conference on Software Testing, Verification and Validation Workshops (ICSTW), 2018, pp. 47-53, doi: 10.1109/ICSTW.2018.00027.'

$ws.Range("A5").Value = "Applying Mutation Testing to Web Applications"
$ws.Range("B5").Value = 'U. Praphamontripong and J. Offutt, "Applying Mutation Testing to Web Applications," 2010 Third International Conference on Software Testing, Verification, and Validation Workshops, 2010, pp. 132-141, doi: 10.1109/ICSTW.2010.38.'

$ws.Range("A6").Value = "An empirical study on the application of mutation testing for a safety-critical industrial software system"
$ws.Range("B6").Value = "Rudolf Ramler, Thomas Wetzlmaier, and Claus Klammer. 2017. An empirical study on the application of mutation testing for a safety-critical industrial software system. In Proceedings of the Symposium on Applied Computing (SAC '17). Association for Computing Machinery, New York, NY, USA, 1401–1408."

# Column A was widened to fit the longer bibliography text.
$ws.Columns.Item(1).ColumnWidth = 94.71

# Active cell moved after the edit.
$ws.Range("A10").Select()
